$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 (the "Bhasha Samarthak Badge" row). This shifts row 7
# ("Konkani Roman") up to become the new row 6, reducing the used range
# from A1:D7 to A1:D6.
$ws.Rows.Item(6).Delete()
